$d = $word.ActiveDocument

# Locate the end of the sentence "... mobile development." which currently
# ends the last diary-entry paragraph (that paragraph also holds two
# trailing manual line-break runs after the text).
$r = $d.Content
$r.Find.Execute("mobile development.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$r.Collapse(0)

# Split the paragraph right after "mobile development." so the sentence
# becomes its own, self-contained paragraph (the trailing line-break runs
# stay behind, now in the next paragraph).
$r.InsertParagraphAfter()
$r.Move(1, 1)

# Build the two new diary entries:
#   (blank line) / 04.11.2020 / (blank line) / <entry text>
# The entry text is typed directly in front of the two manual line breaks
# that used to close the old paragraph, so they now trail the new entry.
# $r currently sits at the start of the blank paragraph created above, so a
# leading carriage return keeps that paragraph empty before "04.11.2020".
$cr = [char]13
$entryText = "While working on the Introduction module I learned that the " + `
  "EditText object method getText will not return a null value in any " + `
  "case. Instead, it returns an empty string if the String inside of the " + `
  "widget is empty."

$newContent = $cr + "04.11.2020" + $cr + $cr + $entryText
$r.InsertAfter($newContent)
